$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 78 (shifts existing rows 78-101 down to 79-102)
$ws.Rows.Item(78).Insert()

# Populate the new row 78 with the new weekly record
$ws.Cells.Item(78, 1).Value = 10
$ws.Cells.Item(78, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(78, 3).Value = "La Araucanía"
$ws.Cells.Item(78, 4).Value = 45027
$ws.Cells.Item(78, 5).Value = 9
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100108
$ws.Cells.Item(78, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(78, 9).Value = 100108004
$ws.Cells.Item(78, 10).Value = "Papaya"
$ws.Cells.Item(78, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(78, 12).Value = "Primera"
$ws.Cells.Item(78, 13).Value = 55
$ws.Cells.Item(78, 14).Value = 3500
$ws.Cells.Item(78, 15).Value = 3500
$ws.Cells.Item(78, 16).Value = 3500
$ws.Cells.Item(78, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(78, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(78, 19).Value = 3500
$ws.Cells.Item(78, 20).Value = 1
